# Complete the first verse of "Replay" in the Cancionero Vigil workbook,
# adding chords/beats for the remaining lyric lines and marking the new
# rows as "played" in column F.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Replay(prueba)")

# Finish the text of the line that was cut short ("dy in my head" -> "...that i ")
$ws.Range("C6").Value = "dy in my head that i "

# Mark the already-existing rows (5-8) as done in the new "played" column F
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = 1

# Fill in the lyrics for rows 7-8 that only had chord/beat data so far
$ws.Range("C7").Value = "cant keep out"
$ws.Range("C8").Value = "got me singing like"

# Add four new rows (9-12) continuing the chord progression E, C, G, D
# (match the centered formatting already used by the D/E columns above)
$ws.Range("D9:E12").HorizontalAlignment = -4108

$ws.Range("C9").Value = "na ra na na "
$ws.Range("D9").Value = "E"
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 1

$ws.Range("C10").Value = "every day"
$ws.Range("D10").Value = "C"
$ws.Range("E10").Value = 6
$ws.Range("F10").Value = 1

$ws.Range("C11").Value = "Its like my ipod "
$ws.Range("D11").Value = "G"
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 1

$ws.Range("C12").Value = "stuck-on replay"
$ws.Range("D12").Value = "D"
$ws.Range("E12").Value = 8
$ws.Range("F12").Value = 1

# Move the active selection to the next empty row, like someone had just
# finished typing the last entry and landed on D13
$ws.Range("D13").Select()
